# Apply the "wrong assumption" diagram fix to both Figure 7g and Figure 7h sheets.
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Figure 7g")
$ws2 = $wb.Worksheets.Item("Figure 7h")

# --- Relabel the AG5 / AG6 header cells -------------------------------------------------
# Clear first so the now-unused shared string is dropped before the new text is
# looked up / appended; this reproduces the exact shared-string slot swap seen in the diff.
$ws1.Range("AG5").ClearContents()
$ws2.Range("AG6").ClearContents()
$ws2.Range("AG6").Value = "conductance fraction"
$ws1.Range("AG5").Value = "current fraction"

# --- Fix the "wrong assumption" formula --------------------------------------------------
# Previously: (AC-AD)/AC  (current difference / conductance loss, decimal parts)
# Now:        AD/AC       (current fraction / conductance fraction)
$ws1.Range("AG7").Formula2 = "=(AD7:AD23)/(AC7:AC23)"
$ws2.Range("AG7").Formula2 = "=(AD7:AD23)/(AC7:AC23)"

# --- Figure 7g: add overall average of the spilled array using ANCHORARRAY -------------
$ws1.Range("AG25").Formula2 = "=AVERAGE(ANCHORARRAY(AG7))"

# --- Figure 7h: add "avvg fraction" label next to the existing average, and highlight ----
$ws2.Range("AF25").Value = "avvg fraction"
$ws2.Range("AF25").Interior.ThemeColor = 6
$ws2.Range("AG25").Interior.ThemeColor = 6

# --- Restore selections to match final state --------------------------------------------
$ws1.Activate()
$ws1.Range("AK15").Select()
$ws2.Activate()
$ws2.Range("AF25").Select()
